# Insert a new data row at row 181 (pushing the existing rows 181:265 down to 182:266)
# and populate it with the new weekly price observation (dated 2021-11-10).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(181).Insert()

$ws.Cells.Item(181, 1).Value = 3
$ws.Cells.Item(181, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(181, 3).Value = "Coquimbo"
$ws.Cells.Item(181, 4).Value = 44510
$ws.Cells.Item(181, 5).Value = 5
$ws.Cells.Item(181, 6).Value = 100112017
$ws.Cells.Item(181, 7).Value = "Apio"
$ws.Cells.Item(181, 8).Value = "Americana (o)"
$ws.Cells.Item(181, 9).Value = "Primera"
$ws.Cells.Item(181, 10).Value = 130
$ws.Cells.Item(181, 11).Value = 9000
$ws.Cells.Item(181, 12).Value = 9000
$ws.Cells.Item(181, 13).Value = 9000
$ws.Cells.Item(181, 14).Value = "$/docena de matas"
$ws.Cells.Item(181, 15).Value = "Pan de Azúcar"
$ws.Cells.Item(181, 16).Value = 1500
$ws.Cells.Item(181, 17).Value = 6
$ws.Cells.Item(181, 18).Value = "Hortaliza"
